# add team page & NorneStore logic
# Adds three new localization keys (wooden_sword / iron_sword / magic_sword)
# as new columns AH, AI, AJ on row 1 (keys), row 2 (Simplified Chinese values)
# for the single "Sheet1" localization table. Row 3 (English) has no values
# for these new keys yet, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - new ID headers
$ws.Range("AH1").Value = "wooden_sword"
$ws.Range("AI1").Value = "iron_sword"
$ws.Range("AJ1").Value = "magic_sword"

# Row 2 - new Simplified Chinese (简体中文) values
$ws.Range("AH2").Value = "木剑"
$ws.Range("AI2").Value = "铁剑"
$ws.Range("AJ2").Value = "神魔一念"

# Move the active selection to the newly added AJ2 cell, and bring the
# new columns into view, matching the author's final view state.
$ws.Range("AJ2").Select()
